$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: K2 no longer carries the hyperlink (it moves to K24 below); its text
# becomes the raw product URL.
# ---------------------------------------------------------------------------
$ws.Range("K2").Hyperlinks.Delete()
$ws.Range("K2").Value = "https://www.mouser.com/ProductDetail/Samsung-Electro-Mechanics/CL05A106MQ5NRNC?qs=xZ%2FP%252Ba9zWqYQV6QCAIWS6w%3D%3D"

# ---------------------------------------------------------------------------
# Row 4: 1uF / U1 capacitor gains Qty/Unit/Total cost + procurement details.
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 15
$ws.Range("B4").Value = 0.12
$ws.Range("C4").Value = 1.8
$ws.Range("B4").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("C4").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("F4").Value = "KEMET"
$ws.Range("G4").Value = "C0402C105M9PACTU"
$ws.Range("H4").Value = "DigiKey"
$ws.Range("I4").Value = "399-C0402C105M9PACTUTR-ND"
$ws.Range("J4").Value = "CAP CER 1UF 6.3V X5R 0402"

# ---------------------------------------------------------------------------
# Row 5: 130pF / U1 capacitor gains Qty/Unit/Total cost + procurement details
# including a manufacturer hyperlink.
# ---------------------------------------------------------------------------
$ws.Range("B5").Value = 0.25
$ws.Range("C5").Value = 1.5
$ws.Range("B5").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("C5").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("G5").Value = "GCG1555G1H131GA01D "
$ws.Range("H5").Value = "Mouser"
$ws.Range("I5").Value = "81-GCG1555G1H131GA1D "
$ws.Range("J5").Value = "Multilayer Ceramic Capacitors MLCC - SMD/SMT 130 pF 50 VDC 2% 0402 X8G AEC-Q200 "
$ws.Range("J5").VerticalAlignment = -4108
$ws.Range("K5").Value = "https://www.mouser.com/ProductDetail/Murata-Electronics/GCG1555G1H131GA01D?qs=QzBtWTOodeUu5id1a%2FuWbQ%3D%3D"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.mouser.com/manufacturer/murataelectronics/", "", "", "Murata Electronics ")
$ws.Range("F5").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Row 23: new DC/DC inverter IC (TC7660EOA) designator row.
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = 1
$ws.Range("D23").Value = "TC7660EOA"
$ws.Range("L23").Value = "Package_SO:SO-8_3.9x4.9mm_P1.27mm"

# ---------------------------------------------------------------------------
# Row 24 (new row): 10uF / U2 capacitor - this is the part that used to be
# referenced from row 2, now moved down with its own hyperlink.
# ---------------------------------------------------------------------------
$ws.Range("A24").Value = 15
$ws.Range("B24").Value = 0.17
$ws.Range("C24").Value = 2.1800000000000002
$ws.Range("B24").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("C24").NumberFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'
$ws.Range("D24").Value = "10uF"
$ws.Range("E24").Value = "U2"
$ws.Range("F24").Value = "Samsung Electro-Mechanics"
$ws.Range("G24").Value = "CL05A106MQ5NRNC"
$ws.Range("H24").Value = "Mouser"
$ws.Range("I24").Value = "187-CL05A106MQ5NRNC "
$ws.Range("J24").Value = "Multilayer Ceramic Capacitors MLCC - SMD/SMT 10uF+/-20% 6.3V X5R 0402"
$ws.Range("K24").Value = "CL05A106MQ5NRNC Samsung Electro-Mechanics | Mouser"
$ws.Range("L24").Value = "Capacitor_SMD:C_0402_1005Metric"
$ws.Hyperlinks.Add($ws.Range("K24"), "https://www.mouser.com/ProductDetail/Samsung-Electro-Mechanics/CL05A106MQ5NRNC?qs=xZ%2FP%252Ba9zWqYQV6QCAIWS6w%3D%3D", "", "", "CL05A106MQ5NRNC Samsung Electro-Mechanics | Mouser")

# ---------------------------------------------------------------------------
# Column width adjustments (reorganized layout).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.5
$ws.Columns.Item(2).ColumnWidth = 13
$ws.Columns.Item(3).ColumnWidth = 13.833333333333334
$ws.Columns.Item(4).ColumnWidth = 22.833333333333332
$ws.Columns.Item(5).ColumnWidth = 10.833333333333334
$ws.Columns.Item(6).ColumnWidth = 28.5
$ws.Columns.Item(7).ColumnWidth = 62.5
$ws.Columns.Item(8).ColumnWidth = 14.5
$ws.Columns.Item(9).ColumnWidth = 29.833333333333332
$ws.Columns.Item(10).ColumnWidth = 78.83333333333333
$ws.Columns.Item(11).ColumnWidth = 110.83333333333333
$ws.Columns.Item(12).ColumnWidth = 70.33333333333333

# ---------------------------------------------------------------------------
# Selection moved to K10 (matches the saved cursor position in the workbook).
# ---------------------------------------------------------------------------
$ws.Range("K10").Select()
